# Weekly refresh of Hortaliza / Agricola del Norte S.A. de Arica - Alcachofa data rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44370
$ws.Range("H2").Value = 'Argentina(o)'
$ws.Range("J2").Value = 140
$ws.Range("K2").Value = 20000
$ws.Range("L2").Value = 21000
$ws.Range("M2").Value = 20429
$ws.Range("N2").Value = '$/caja 50 unidades'
$ws.Range("P2").Value = 409
$ws.Range("Q2").Value = 50
# Row 3
$ws.Range("D3").Value = 44370
$ws.Range("J3").Value = 80
$ws.Range("K3").Value = 22000
$ws.Range("L3").Value = 23000
$ws.Range("M3").Value = 22500
$ws.Range("P3").Value = 562
# Row 4
$ws.Range("D4").Value = 44377
$ws.Range("J4").Value = 150
$ws.Range("K4").Value = 20000
$ws.Range("L4").Value = 21000
$ws.Range("M4").Value = 20333
$ws.Range("P4").Value = 508
# Row 5
$ws.Range("D5").Value = 44377
$ws.Range("H5").Value = 'Symphony'
$ws.Range("I5").Value = 'Primera'
$ws.Range("J5").Value = 60
$ws.Range("K5").Value = 21000
$ws.Range("L5").Value = 22000
$ws.Range("M5").Value = 21500
$ws.Range("N5").Value = '$/caja 40 unidades'
$ws.Range("P5").Value = 538
$ws.Range("Q5").Value = 40
# Row 6
$ws.Range("D6").Value = 44435
$ws.Range("H6").Value = 'Madrigal'
$ws.Range("J6").Value = 160
$ws.Range("K6").Value = 19000
$ws.Range("L6").Value = 20000
$ws.Range("M6").Value = 19500
$ws.Range("P6").Value = 488
# Row 7
$ws.Range("D7").Value = 44167
$ws.Range("H7").Value = 'Española'
$ws.Range("J7").Value = 160
$ws.Range("K7").Value = 13000
$ws.Range("L7").Value = 14000
$ws.Range("M7").Value = 13500
$ws.Range("N7").Value = '$/caja 30 unidades'
$ws.Range("O7").Value = 'Región Metropolitana'
$ws.Range("P7").Value = 450
$ws.Range("Q7").Value = 30
# Row 8
$ws.Range("D8").Value = 44468
$ws.Range("H8").Value = 'Argentina(o)'
$ws.Range("J8").Value = 120
$ws.Range("K8").Value = 17000
$ws.Range("L8").Value = 18000
$ws.Range("M8").Value = 17500
$ws.Range("N8").Value = '$/caja 50 unidades'
$ws.Range("P8").Value = 350
$ws.Range("Q8").Value = 50
# Row 9
$ws.Range("D9").Value = 44363
$ws.Range("J9").Value = 160
$ws.Range("K9").Value = 19000
$ws.Range("L9").Value = 20000
$ws.Range("M9").Value = 19500
$ws.Range("P9").Value = 488
# Row 10
$ws.Range("D10").Value = 44356
$ws.Range("K10").Value = 19000
$ws.Range("L10").Value = 20000
$ws.Range("M10").Value = 19500
$ws.Range("P10").Value = 390
# Row 11
$ws.Range("D11").Value = 44483
$ws.Range("J11").Value = 120
$ws.Range("K11").Value = 14000
$ws.Range("L11").Value = 15000
$ws.Range("M11").Value = 14500
$ws.Range("P11").Value = 362
# Row 12
$ws.Range("D12").Value = 44391
$ws.Range("H12").Value = 'Madrigal'
$ws.Range("J12").Value = 140
# Row 13
$ws.Range("D13").Value = 44398
$ws.Range("J13").Value = 170
$ws.Range("K13").Value = 21000
$ws.Range("L13").Value = 22000
$ws.Range("M13").Value = 21500
$ws.Range("P13").Value = 538
# Row 14
$ws.Range("D14").Value = 44384
$ws.Range("J14").Value = 80
$ws.Range("K14").Value = 21000
$ws.Range("L14").Value = 22000
$ws.Range("M14").Value = 21500
$ws.Range("P14").Value = 538
# Row 15
$ws.Range("D15").Value = 44384
$ws.Range("I15").Value = 'Segunda'
$ws.Range("J15").Value = 30
$ws.Range("K15").Value = 19000
$ws.Range("L15").Value = 20000
$ws.Range("M15").Value = 19333
$ws.Range("N15").Value = '$/caja 50 unidades'
$ws.Range("P15").Value = 387
$ws.Range("Q15").Value = 50
# Row 16
$ws.Range("D16").Value = 44384
$ws.Range("H16").Value = 'Symphony'
$ws.Range("J16").Value = 100
$ws.Range("K16").Value = 20000
$ws.Range("L16").Value = 21000
$ws.Range("M16").Value = 20400
$ws.Range("P16").Value = 510
# Row 17
$ws.Range("D17").Value = 44412
$ws.Range("H17").Value = 'Symphony'
$ws.Range("J17").Value = 240
$ws.Range("K17").Value = 21000
$ws.Range("L17").Value = 22000
$ws.Range("M17").Value = 21500
$ws.Range("P17").Value = 538
# Row 18
$ws.Range("D18").Value = 44405
$ws.Range("J18").Value = 200
# Row 19
$ws.Range("D19").Value = 44426
$ws.Range("H19").Value = 'Madrigal'
$ws.Range("J19").Value = 150
$ws.Range("K19").Value = 19000
$ws.Range("L19").Value = 20000
$ws.Range("M19").Value = 19500
$ws.Range("N19").Value = '$/caja 40 unidades'
$ws.Range("P19").Value = 488
$ws.Range("Q19").Value = 40
# Row 20
$ws.Range("D20").Value = 44160
$ws.Range("J20").Value = 160
$ws.Range("K20").Value = 14000
$ws.Range("L20").Value = 15000
$ws.Range("M20").Value = 14500
$ws.Range("P20").Value = 362
# Row 21
$ws.Range("D21").Value = 44482
$ws.Range("H21").Value = 'Madrigal'
$ws.Range("J21").Value = 200
$ws.Range("K21").Value = 14000
$ws.Range("L21").Value = 15000
$ws.Range("M21").Value = 14500
$ws.Range("N21").Value = '$/caja 40 unidades'
$ws.Range("O21").Value = 'Región de Coquimbo'
$ws.Range("P21").Value = 362
$ws.Range("Q21").Value = 40
# Row 22
$ws.Range("D22").Value = 44433
$ws.Range("H22").Value = 'Madrigal'
$ws.Range("J22").Value = 160
$ws.Range("K22").Value = 19000
$ws.Range("L22").Value = 20000
$ws.Range("M22").Value = 19500
$ws.Range("N22").Value = '$/caja 40 unidades'
$ws.Range("P22").Value = 488
$ws.Range("Q22").Value = 40
# Row 24
$ws.Range("D24").Value = 44419
$ws.Range("H24").Value = 'Symphony'
$ws.Range("J24").Value = 150
$ws.Range("K24").Value = 21000
$ws.Range("L24").Value = 22000
$ws.Range("M24").Value = 21500
$ws.Range("N24").Value = '$/caja 50 unidades'
$ws.Range("P24").Value = 430
$ws.Range("Q24").Value = 50
